# Updated symbol list on Thu Jan 26 20:35:29 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) quotes that changed
# since the last data pull, keeping values stored as text (matching the
# existing inline-string / percent-as-text convention used in this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    # Leading apostrophe forces Excel to store the value as literal text
    # (so numbers like '306.71' and percentages like '1.07%' are kept as
    # strings instead of being parsed into numeric/percentage values).
    $cell.Value = "'" + $text
    # Re-apply the Normal style so no stray number-format/quote-prefix
    # style id gets attached to the cell.
    $cell.Style = "Normal"
}

Set-TextCell $ws 2 4 "306.71"
Set-TextCell $ws 2 5 "1.07%"
Set-TextCell $ws 3 4 "35.97"
Set-TextCell $ws 3 5 "0.73%"
Set-TextCell $ws 4 4 "5.018"
Set-TextCell $ws 4 5 "-1.03%"
Set-TextCell $ws 5 4 "0.08095"
Set-TextCell $ws 5 5 "0.34%"
Set-TextCell $ws 6 4 "1.920"
Set-TextCell $ws 6 5 "-0.74%"
Set-TextCell $ws 7 5 "2.47%"
Set-TextCell $ws 8 4 "7.863"
Set-TextCell $ws 8 5 "0.79%"
Set-TextCell $ws 9 4 "0.9304"
Set-TextCell $ws 9 5 "0.10%"
Set-TextCell $ws 10 4 "0.1252"
Set-TextCell $ws 10 5 "-18.72%"
Set-TextCell $ws 11 4 "0.1910"
Set-TextCell $ws 11 5 "0.67%"
Set-TextCell $ws 12 4 "0.09256"
Set-TextCell $ws 12 5 "2.67%"
Set-TextCell $ws 13 4 "0.03504"
Set-TextCell $ws 13 5 "1.62%"
Set-TextCell $ws 14 4 "0.09921"
Set-TextCell $ws 14 5 "0.70%"
Set-TextCell $ws 15 4 "0.001424"
Set-TextCell $ws 15 5 "1.98%"
Set-TextCell $ws 16 4 "0.006368"
Set-TextCell $ws 16 5 "10.43%"
Set-TextCell $ws 17 4 "3.609"
Set-TextCell $ws 17 5 "1.83%"
Set-TextCell $ws 18 4 "3.005"
Set-TextCell $ws 18 5 "0.73%"
Set-TextCell $ws 19 4 "0.3438"
Set-TextCell $ws 19 5 "-0.11%"
Set-TextCell $ws 20 4 "5.175"
Set-TextCell $ws 20 5 "2.99%"
Set-TextCell $ws 21 4 "0.1294"
Set-TextCell $ws 21 5 "-0.32%"
Set-TextCell $ws 22 5 "1.60%"
Set-TextCell $ws 23 4 "0.04411"
Set-TextCell $ws 23 5 "-1.68%"
Set-TextCell $ws 24 4 "0.001233"
Set-TextCell $ws 24 5 "1.88%"
Set-TextCell $ws 25 4 "0.004731"
Set-TextCell $ws 25 5 "-1.65%"
Set-TextCell $ws 26 5 "5.99%"
Set-TextCell $ws 27 4 "0.0003127"
Set-TextCell $ws 27 5 "3.73%"
Set-TextCell $ws 39 4 "0.01962"
Set-TextCell $ws 39 5 "4.80%"
Set-TextCell $ws 40 4 "0.05244"
Set-TextCell $ws 40 5 "9.30%"
Set-TextCell $ws 41 4 "0.007558"
Set-TextCell $ws 41 5 "3.22%"
Set-TextCell $ws 42 4 "0.01018"
Set-TextCell $ws 42 5 "-3.86%"
Set-TextCell $ws 43 4 "0.1374"
Set-TextCell $ws 43 5 "2.12%"
Set-TextCell $ws 44 4 "0.002100"
Set-TextCell $ws 44 5 "-0.19%"
Set-TextCell $ws 45 4 "0.01068"
Set-TextCell $ws 45 5 "9.74%"
Set-TextCell $ws 46 4 "0.00006362"
Set-TextCell $ws 46 5 "2.10%"
Set-TextCell $ws 47 5 "0.26%"
Set-TextCell $ws 48 4 "63.57"
Set-TextCell $ws 48 5 "-1.70%"
Set-TextCell $ws 49 4 "0.001658"
Set-TextCell $ws 49 5 "-0.07%"
Set-TextCell $ws 50 4 "0.00002099"
Set-TextCell $ws 50 5 "0.26%"
Set-TextCell $ws 51 4 "0.0001999"
Set-TextCell $ws 51 5 "0.26%"
